$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2: CART ----
$ws.Range("B2").Value = "{'max_depth': 5, 'min_samples_leaf': 5, 'min_samples_split': 10}"
$ws.Range("C2").Value = 0.7293612987320773
$ws.Range("D2").Value = "DecisionTreeClassifier(max_depth=5, min_samples_leaf=5, min_samples_split=10)"
$ws.Range("M2").Value = 16.74492311477661

# ---- Row 3: Random Forest ----
$ws.Range("B3").Value = "{'max_depth': 20, 'min_samples_split': 10, 'n_estimators': 150}"
$ws.Range("C3").Value = 0.7551735987676265
$ws.Range("D3").Value = "RandomForestClassifier(max_depth=20, min_samples_split=10, n_estimators=150)"
$ws.Range("E3").Value = "[[181  96]`n [ 86 363]]"
$ws.Range("F3").Value = 363
$ws.Range("G3").Value = 96
$ws.Range("H3").Value = 86
$ws.Range("I3").Value = 181
$ws.Range("J3").Value = 0.7477555502550995
$ws.Range("K3").Value = 0.7493112947658402
$ws.Range("L3").Value = 0.7483876150489321
$ws.Range("M3").Value = 625.2701163291931

# ---- Row 4: LightGBM ----
$ws.Range("B4").Value = "{'learning_rate': 0.15, 'n_estimators': 50, 'num_leaves': 50}"
$ws.Range("C4").Value = 0.7613757554212583
$ws.Range("D4").Value = "LGBMClassifier(learning_rate=0.15, n_estimators=50, num_leaves=50)"
$ws.Range("E4").Value = "[[183  94]`n [ 90 359]]"
$ws.Range("F4").Value = 359
$ws.Range("G4").Value = 94
$ws.Range("H4").Value = 90
$ws.Range("I4").Value = 183
$ws.Range("J4").Value = 0.7458833894715893
$ws.Range("K4").Value = 0.7465564738292011
$ws.Range("L4").Value = 0.7461968200448346
$ws.Range("M4").Value = 222.5531287193298

# ---- Row 5: XGBoost ----
$ws.Range("B5").Value = "{'learning_rate': 0.05, 'max_depth': 3, 'n_estimators': 150}"
$ws.Range("C5").Value = 0.7692913852352175
$ws.Range("D5").Value = "XGBClassifier(base_score=None, booster=None, callbacks=None,`n              colsample_bylevel=None, colsample_bynode=None,`n              colsample_bytree=None, device=None, early_stopping_rounds=None,`n              enable_categorical=True, eval_metric=None, feature_types=None,`n              gamma=None, grow_policy=None, importance_type=None,`n              interaction_constraints=None, learning_rate=0.05, max_bin=None,`n              max_cat_threshold=None, max_cat_to_onehot=None,`n              max_delta_step=None, max_depth=3, max_leaves=None,`n              min_child_weight=None, missing=nan, monotone_constraints=None,`n              multi_strategy=None, n_estimators=150, n_jobs=None,`n              num_parallel_tree=None, random_state=None, ...)"
$ws.Range("E5").Value = "[[186  91]`n [ 75 374]]"
$ws.Range("F5").Value = 374
$ws.Range("G5").Value = 91
$ws.Range("H5").Value = 75
$ws.Range("I5").Value = 186
$ws.Range("J5").Value = 0.7693298645265477
$ws.Range("K5").Value = 0.7713498622589532
$ws.Range("L5").Value = 0.7699511743266545
$ws.Range("M5").Value = 274.0794744491577

# Re-assigning the multi-line confusion-matrix strings can trigger an
# automatic custom-row-height (wrap) side effect; AutoFit restores the
# rows to their natural (default) height so the output matches the
# original formatting.
$ws.Rows("2:5").AutoFit()
